$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A - this shifts the existing
#    FirstName/LastName/Email columns to B/C/D and preserves their widths.
$ws.Range("A:A").Insert()

# 2. Give the new column A a width (closest achievable to 14.140625 chars).
$ws.Columns.Item(1).ColumnWidth = 13.25

# 3. Fill in the new "Start" / "TestDATA" / "END" column.
$ws.Range("A1").Value = "Start"
$ws.Range("A2").Value = "TestDATA"
$ws.Range("A3").Value = "END"

# 4. The hyperlink that used to live on C2 now points at the right cell but
#    its stored range reference is stale after the column insert - drop it
#    and re-create it on D2 (the Email cell's new location).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:rahulgautamvvs@gmail.com") | Out-Null

# 5. Formatting: header row + column A get a yellow fill and a thin border;
#    the rest of the data region gets just the thin border; the hyperlink
#    cell keeps the Hyperlink style plus the border.
$ws.Range("A1:D1").Interior.Color = 65535
$ws.Range("A1:D1").Borders.LineStyle = 1

$ws.Range("A2:A3").Interior.Color = 65535
$ws.Range("A2:A3").Borders.LineStyle = 1

$ws.Range("B2:D3").Borders.LineStyle = 1

$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D2").Borders.LineStyle = 1

# 6. Move the active selection, matching the saved view state.
$ws.Range("D9").Select()

# 7. Page setup orientation, as recorded in the saved file.
$ws.PageSetup.Orientation = 1
